$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 3368019.8
$ws.Range("I33").Value = 7405585.5
$ws.Range("J33").Value = 3381.5
$ws.Range("K33").Value = 7405585.5
$ws.Range("L33").Value = 3381.5
$ws.Range("M33").Value = -7405356.5
$ws.Range("N33").Value = -3839.5
$ws.Range("H75").Value = 18642.5
$ws.Range("I75").Value = 18642.5
$ws.Range("K75").Value = 18642.5
$ws.Range("M75").Value = -17706.5
$ws.Range("H78").Value = 18642.5
$ws.Range("I78").Value = 18642.5
$ws.Range("K78").Value = 55927.5
$ws.Range("M78").Value = -51247.5
$ws.Range("H92").Value = 515.3226
$ws.Range("I92").Value = 517.0741
$ws.Range("K92").Value = 517.0741
$ws.Range("M92").Value = 730.9259
$ws.Range("H99").Value = 465
$ws.Range("I99").Value = 496.2
$ws.Range("J99").Value = 426
$ws.Range("K99").Value = 1488.6
$ws.Range("L99").Value = 1278
$ws.Range("M99").Value = 9.400000000000091
$ws.Range("N99").Value = -4274
$ws.Range("H101").Value = 1114.4286
$ws.Range("I101").Value = 1374.75
$ws.Range("J101").Value = 767.3333
$ws.Range("K101").Value = 4124.25
$ws.Range("L101").Value = 2301.9999
$ws.Range("M101").Value = -2502.25
$ws.Range("N101").Value = -5545.9999
$ws.Range("H103").Value = 391
$ws.Range("I103").Value = 359
$ws.Range("K103").Value = 1077
$ws.Range("M103").Value = -491
$ws.Range("H104").Value = 171.33333
$ws.Range("I104").Value = 171.33333
$ws.Range("K104").Value = 513.99999
$ws.Range("M104").Value = 1233.00001
$ws.Range("H137").Value = 16666.408
$ws.Range("I137").Value = 7987.4443
$ws.Range("J137").Value = 22674.924
$ws.Range("K137").Value = 23962.3329
$ws.Range("L137").Value = 68024.772
$ws.Range("M137").Value = -21412.3329
$ws.Range("N137").Value = -73124.772
$ws.Range("H138").Value = 1769.6666
$ws.Range("I138").Value = 1495
$ws.Range("J138").Value = 4173
$ws.Range("K138").Value = 4485
$ws.Range("L138").Value = 12519
$ws.Range("M138").Value = 655
$ws.Range("N138").Value = -22799
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6327.659
$ws.Range("I32").Value = 5571.1465
$ws.Range("J32").Value = 16666.666
$ws.Range("K32").Value = 5571.1465
$ws.Range("L32").Value = 16666.666
$ws.Range("M32").Value = -5284.1465
$ws.Range("N32").Value = -17240.666
$ws.Range("H45").Value = 6447.5654
$ws.Range("J45").Value = 1919.6
$ws.Range("L45").Value = 1919.6
$ws.Range("N45").Value = -2673.6
$ws.Range("H61").Value = 19733.285
$ws.Range("I61").Value = 16906.375
$ws.Range("J61").Value = 28779.4
$ws.Range("K61").Value = 16906.375
$ws.Range("L61").Value = 28779.4
$ws.Range("M61").Value = -16694.375
$ws.Range("N61").Value = -29203.4
$ws.Range("H132").Value = 981.6857
$ws.Range("I132").Value = 932.4828
$ws.Range("J132").Value = 1219.5
$ws.Range("K132").Value = 2797.4484
$ws.Range("L132").Value = 3658.5
$ws.Range("M132").Value = -267.4484000000002
$ws.Range("N132").Value = -8718.5
$ws.Range("H136").Value = 19733.285
$ws.Range("I136").Value = 16906.375
$ws.Range("J136").Value = 28779.4
$ws.Range("K136").Value = 50719.125
$ws.Range("L136").Value = 86338.20000000001
$ws.Range("M136").Value = -48169.125
$ws.Range("N136").Value = -91438.20000000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 5064.65
$ws.Range("I99").Value = 5146.972
$ws.Range("J99").Value = 4323.75
$ws.Range("K99").Value = 5146.972
$ws.Range("L99").Value = 4323.75
$ws.Range("M99").Value = -3648.972
$ws.Range("N99").Value = -7319.75
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3435.261
$ws.Range("I31").Value = 1751.7778
$ws.Range("J31").Value = 4517.5
$ws.Range("K31").Value = 1751.7778
$ws.Range("L31").Value = 4517.5
$ws.Range("M31").Value = -1456.7778
$ws.Range("N31").Value = -5107.5
$ws.Range("H34").Value = 3435.261
$ws.Range("I34").Value = 1751.7778
$ws.Range("J34").Value = 4517.5
$ws.Range("K34").Value = 1751.7778
$ws.Range("L34").Value = 4517.5
$ws.Range("M34").Value = -1549.7778
$ws.Range("N34").Value = -4921.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 4885.467
$ws.Range("I131").Value = 1709.9
$ws.Range("J131").Value = 5792.7715
$ws.Range("K131").Value = 5129.700000000001
$ws.Range("L131").Value = 17378.3145
$ws.Range("M131").Value = -89.70000000000073
$ws.Range("N131").Value = -27458.3145
$ws.Range("H140").Value = 1617.8334
$ws.Range("I140").Value = 1617.8334
$ws.Range("K140").Value = 4853.5002
$ws.Range("M140").Value = 326.4997999999996
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 27987.334
$ws.Range("J136").Value = 27987.334
$ws.Range("L136").Value = 83962.00199999999
$ws.Range("N136").Value = -89062.00199999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3438.8572
$ws.Range("I68").Value = 3028
$ws.Range("J68").Value = 3812.3635
$ws.Range("K68").Value = 3028
$ws.Range("L68").Value = 3812.3635
$ws.Range("M68").Value = -2279
$ws.Range("N68").Value = -5310.363499999999
$ws.Range("H71").Value = 3438.8572
$ws.Range("I71").Value = 3028
$ws.Range("J71").Value = 3812.3635
$ws.Range("K71").Value = 15140
$ws.Range("L71").Value = 19061.8175
$ws.Range("M71").Value = -11396
$ws.Range("N71").Value = -26549.8175
$ws.Range("H99").Value = 29701.6
$ws.Range("I99").Value = 29701.6
$ws.Range("K99").Value = 29701.6
$ws.Range("M99").Value = -26706.6
$ws.Range("H122").Value = 3769.6553
$ws.Range("I122").Value = 2894.5908
$ws.Range("J122").Value = 6519.857
$ws.Range("K122").Value = 8683.7724
$ws.Range("L122").Value = 19559.571
$ws.Range("M122").Value = -6233.7724
$ws.Range("N122").Value = -24459.571
$ws.Range("H132").Value = 7737.5713
$ws.Range("I132").Value = 7485.091
$ws.Range("J132").Value = 8663.333000000001
$ws.Range("K132").Value = 22455.273
$ws.Range("L132").Value = 25989.999
$ws.Range("M132").Value = -19925.273
$ws.Range("N132").Value = -31049.999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 22148
$ws.Range("I52").Value = 22148
$ws.Range("K52").Value = 22148
$ws.Range("M52").Value = -21922
$ws.Range("H132").Value = 160213.27
$ws.Range("I132").Value = 237078
$ws.Range("J132").Value = 30134.5
$ws.Range("K132").Value = 711234
$ws.Range("L132").Value = 90403.5
$ws.Range("M132").Value = -708704
$ws.Range("N132").Value = -95463.5
